$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.148.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.552.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.553.95'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.51%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("E10").Value = '  -2.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.409'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.163.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000202'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.553.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.29%  '

$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.188.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '421.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.604'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.700.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.36%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.27%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.555.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("E33").Value = '  +3.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.81%  '

$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.27%  '

$ws.Range("E39").Value = '  -6.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0834'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.871'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.67%  '

$ws.Range("E45").Value = '  -5.83%  '

$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.21%  '
